# Insert a new data row at row 125 ("Fruta, Feria Lagunitas de Puerto Montt - Uva"
# sheet), pushing the existing rows 125-218 down to 126-219, and fill the newly
# inserted row with the new "Rosada pastilla" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 125 - this shifts rows 125..218
# down to 126..219 and extends the sheet dimension to A1:T219.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new record.
$ws.Range("A125").Value = 4
$ws.Range("B125").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C125").Value = "Los Lagos"
$ws.Range("D125").Value = 44651
$ws.Range("E125").Value = 10
$ws.Range("F125").Value = "Fruta"
$ws.Range("G125").Value = 100109
$ws.Range("H125").Value = "Uva"
$ws.Range("I125").Value = 100109001
$ws.Range("J125").Value = "Uva"
$ws.Range("K125").Value = "Rosada pastilla"
$ws.Range("L125").Value = "Primera"
$ws.Range("M125").Value = 300
$ws.Range("N125").Value = 14000
$ws.Range("O125").Value = 15000
$ws.Range("P125").Value = 14500
$ws.Range("Q125").Value = "$/bandeja 10 kilos"
$ws.Range("R125").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S125").Value = 1450
$ws.Range("T125").Value = 10
